$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.866.02"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "1.620.24"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.250"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0617"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "1.844.46"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").Value = "1.609.72"
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.93%  "
$ws.Range("D16").Value = "25.878.43"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.43%  "
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("E24").Value = "  +3.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0478"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("E32").Value = "  -3.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.51%  "
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("D36").Value = "1.126.76"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.841"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.08%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.514"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.48%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "1.754.94"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.749"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.75%  "
$ws.Range("E44").Value = "  -5.07%  "
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "54.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.27%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.14%  "
